$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update status text "Ready for handoff" -> "In Translation" on every
#    sheet/cell where it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Narrow the date-time columns from ~17.22 characters to ~13.41 characters.
#    (Overview columns E & F, zh-cn column C, de-de column C.)
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
